# Trade #116 closed at 2026-02-17 16:03:25 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1198.86
$wsSummary.Range("B4").Value = -1.15
$wsSummary.Range("B6").Value = 116
$wsSummary.Range("B7").Value = 42
$wsSummary.Range("B9").Value = 36.21

# ---- Strategy Status sheet ----
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 98.86
$wsStrategy.Range("D4").Value = 116
$wsStrategy.Range("E4").Value = -1.15
$wsStrategy.Range("F4").Value = -1.14
$wsStrategy.Range("G4").Value = 36.21

# ---- Append new trade row (#116) to both "All Trades" and "MarketMaking" sheets ----
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRowIndex = 117

    $ws.Cells.Item($newRowIndex, 1).Value = 116
    # Force column B to text so the date-like string "2026-02-17" is not
    # auto-converted into a date serial number (matches source inlineStr).
    $ws.Cells.Item($newRowIndex, 2).NumberFormat = "@"
    $ws.Cells.Item($newRowIndex, 2).Value = "2026-02-17"
    $ws.Cells.Item($newRowIndex, 3).Value = "16:03:18"
    $ws.Cells.Item($newRowIndex, 4).Value = "MarketMaking"
    $ws.Cells.Item($newRowIndex, 5).Value = "UP"
    $ws.Cells.Item($newRowIndex, 6).Value = 0.86
    $ws.Cells.Item($newRowIndex, 7).Value = 0.87
    $ws.Cells.Item($newRowIndex, 8).Value = "CLOSED"
    $ws.Cells.Item($newRowIndex, 9).Value = 1.1628
    $ws.Cells.Item($newRowIndex, 10).Value = 0.01
    $ws.Cells.Item($newRowIndex, 11).Value = 98.86
    $ws.Cells.Item($newRowIndex, 12).Value = 0
    $ws.Cells.Item($newRowIndex, 13).Value = 0
    $ws.Cells.Item($newRowIndex, 14).Value = 0.6
    $ws.Cells.Item($newRowIndex, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($newRowIndex, 16).Value = "early_exit"
    $ws.Cells.Item($newRowIndex, 17).Value = 0.14
}
